# The sheet lists insurer names in column B and their prices in column C.
# This edit:
#   - adds a new top entry "Niva Bupa (formerly known as Max Bupa)" / "₹1,396"
#   - updates "Oriental"'s price from ₹1,344 to ₹1,460
#   - keeps "National Insurance" / "₹1,503" as-is
#   - drops the former "Bajaj Allianz" / "₹1,672" row
# Net effect: row 2 -> Niva Bupa, row 3 -> Oriental (new price), row 4 -> National Insurance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Niva Bupa (formerly known as Max Bupa)"
$ws.Range("C2").Value = "₹1,396"
$ws.Range("B3").Value = "Oriental"
$ws.Range("C3").Value = "₹1,460"
$ws.Range("B4").Value = "National Insurance"
$ws.Range("C4").Value = "₹1,503"
